$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 134.75
$ws.Range("I33").Value = 113
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 113
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 116
$ws.Range("N33").Value = -658
$ws.Range("H40").Value = 6723.7
$ws.Range("J40").Value = 7132.5
$ws.Range("L40").Value = 7132.5
$ws.Range("N40").Value = -7482.5
$ws.Range("H104").Value = 908.5
$ws.Range("I104").Value = 908.5
$ws.Range("K104").Value = 2725.5
$ws.Range("M104").Value = -978.5
$ws.Range("H105").Value = 17629.875
$ws.Range("J105").Value = 17629.875
$ws.Range("L105").Value = 17629.875
$ws.Range("N105").Value = -24617.875
$ws.Range("H112").Value = 2724.077
$ws.Range("J112").Value = 2591.3
$ws.Range("L112").Value = 7773.900000000001
$ws.Range("N112").Value = -9989.900000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5884246.5
$ws.Range("I32").Value = 2012
$ws.Range("K32").Value = 2012
$ws.Range("M32").Value = -1725
$ws.Range("H45").Value = 3275.5334
$ws.Range("I45").Value = 2193.889
$ws.Range("K45").Value = 2193.889
$ws.Range("M45").Value = -1816.889
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H61").Value = 2750
$ws.Range("I61").Value = 2750
$ws.Range("K61").Value = 2750
$ws.Range("M61").Value = -2538
$ws.Range("H88").Value = 1617.7273
$ws.Range("I88").Value = 2158.3333
$ws.Range("K88").Value = 2158.3333
$ws.Range("M88").Value = -1752.3333
$ws.Range("H91").Value = 1617.7273
$ws.Range("I91").Value = 2158.3333
$ws.Range("K91").Value = 2158.3333
$ws.Range("M91").Value = -754.3332999999998
$ws.Range("H105").Value = 24998.334
$ws.Range("J105").Value = 24998.334
$ws.Range("L105").Value = 24998.334
$ws.Range("N105").Value = -31986.334
$ws.Range("H106").Value = 32599.8
$ws.Range("J106").Value = 32599.8
$ws.Range("L106").Value = 32599.8
$ws.Range("N106").Value = -35123.8
$ws.Range("H132").Value = 1996.3334
$ws.Range("I132").Value = 1996.3334
$ws.Range("K132").Value = 5989.0002
$ws.Range("M132").Value = -3459.0002
$ws.Range("H136").Value = 2750
$ws.Range("I136").Value = 2750
$ws.Range("K136").Value = 8250
$ws.Range("M136").Value = -5700

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 886.5
$ws.Range("I36").Value = 886.5
$ws.Range("K36").Value = 886.5
$ws.Range("M36").Value = -352.5
$ws.Range("H54").Value = 6511.5
$ws.Range("J54").Value = 6798
$ws.Range("L54").Value = 6798
$ws.Range("N54").Value = -7766

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1483.3334
$ws.Range("I16").Value = 1483.3334
$ws.Range("K16").Value = 1483.3334
$ws.Range("M16").Value = -1196.3334
$ws.Range("H31").Value = 5120.8965
$ws.Range("I31").Value = 2227.818
$ws.Range("K31").Value = 2227.818
$ws.Range("M31").Value = -1932.818
$ws.Range("H34").Value = 5120.8965
$ws.Range("I34").Value = 2227.818
$ws.Range("K34").Value = 2227.818
$ws.Range("M34").Value = -2025.818
$ws.Range("H113").Value = 1483.3334
$ws.Range("I113").Value = 1483.3334
$ws.Range("K113").Value = 1483.3334
$ws.Range("M113").Value = 686.6666
$ws.Range("H134").Value = 3598.2
$ws.Range("I134").Value = 3598.2
$ws.Range("K134").Value = 10794.6
$ws.Range("M134").Value = -8259.599999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 169.8
$ws.Range("J15").Value = 187.5
$ws.Range("L15").Value = 562.5
$ws.Range("N15").Value = -842.5
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H59").Value = 950
$ws.Range("I59").Value = 950
$ws.Range("K59").Value = 2850
$ws.Range("M59").Value = -2310
$ws.Range("H92").Value = 2128.3845
$ws.Range("I92").Value = 925
$ws.Range("J92").Value = 2663.2222
$ws.Range("K92").Value = 2775
$ws.Range("L92").Value = 7989.6666
$ws.Range("M92").Value = -1527
$ws.Range("N92").Value = -10485.6666
$ws.Range("H107").Value = 476.7097
$ws.Range("I107").Value = 218.53334
$ws.Range("J107").Value = 718.75
$ws.Range("K107").Value = 655.6000200000001
$ws.Range("L107").Value = 2156.25
$ws.Range("M107").Value = 1264.39998
$ws.Range("N107").Value = -5996.25
$ws.Range("H117").Value = 1582
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("H131").Value = 1865.8572
$ws.Range("I131").Value = 1274.8334
$ws.Range("K131").Value = 3824.5002
$ws.Range("M131").Value = 1215.4998
$ws.Range("H132").Value = 4783.5713
$ws.Range("I132").Value = 4745.3335
$ws.Range("J132").Value = 4812.25
$ws.Range("K132").Value = 42708.0015
$ws.Range("L132").Value = 43310.25
$ws.Range("M132").Value = -40178.0015
$ws.Range("N132").Value = -48370.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2250.6924
$ws.Range("I132").Value = 2250.6924
$ws.Range("K132").Value = 6752.0772
$ws.Range("M132").Value = -4222.0772

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1423.3334
$ws.Range("I55").Value = 1263
$ws.Range("K55").Value = 1263
$ws.Range("M55").Value = -1090
$ws.Range("H106").Value = 8414
$ws.Range("J106").Value = 8414
$ws.Range("L106").Value = 8414
$ws.Range("N106").Value = -10938

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 8000
$ws.Range("I64").Value = 8000
$ws.Range("K64").Value = 8000
$ws.Range("M64").Value = -7752
$ws.Range("H67").Value = 8000
$ws.Range("I67").Value = 8000
$ws.Range("K67").Value = 8000
$ws.Range("M67").Value = -7142
$ws.Range("H100").Value = 2301
$ws.Range("I100").Value = 2401.3333
$ws.Range("K100").Value = 4802.6666
$ws.Range("M100").Value = -4261.6666
$ws.Range("H105").Value = 32224
$ws.Range("J105").Value = 32224
$ws.Range("L105").Value = 32224
$ws.Range("N105").Value = -39212
